$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.029.92'
$ws.Range('E2').Value = '  -0.03%  '
$ws.Range('D3').Value = '3.414.42'
$ws.Range('E3').Value = '  +0.10%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '408.88'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.39%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '129.20'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -3.45%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.637'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +7.07%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('E9').Value = '  +7.22%  '
$ws.Range('E10').Value = '  +19.41%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '42.42'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.76%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000225'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +73.63%  '
$ws.Range('D14').Value = '3.961.38'
$ws.Range('E14').Value = '  +0.09%  '
$ws.Range('E15').Value = '  +5.80%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '20.77'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +4.51%  '
$ws.Range('D17').Value = '3.414.31'
$ws.Range('E17').Value = '  -0.61%  '
$ws.Range('E18').Value = '  +10.25%  '
$ws.Range('E19').Value = '  +5.58%  '
$ws.Range('D20').Value = '61.961.41'
$ws.Range('E20').Value = '  -0.25%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '408.88'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +30.46%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '89.36'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +6.21%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.17'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.58%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.06'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.47%  '
$ws.Range('E25').Value = '  +2.30%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '33.08'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +11.82%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.92'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +9.00%  '
$ws.Range('E28').Value = '  -0.51%  '
$ws.Range('E29').Value = '  +0.03%  '
$ws.Range('E30').Value = '  -2.04%  '
$ws.Range('B31').Value = 'Cosmos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '11.87'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +4.19%  '
$ws.Range('B32').Value = 'Kaspa'
$ws.Range('C32').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.171'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -1.79%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.115'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.01%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '42.67'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.22%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0498'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +3.15%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '54.00'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +4.18%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.999'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.05%  '
$ws.Range('B39').Value = 'LidoDAOToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.35'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.95%  '
$ws.Range('B40').Value = 'Stellar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.134'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +6.86%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.91'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -1.03%  '
$ws.Range('E42').Value = '  +3.93%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '141.67'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +3.20%  '
$ws.Range('E44').Value = '  -0.99%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.11'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +1.59%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.41'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +8.40%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '16.64'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.69%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '21.83'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +2.80%  '
$ws.Range('D49').Value = '2.112.14'
$ws.Range('E49').Value = '  -0.44%  '
$ws.Range('E50').Value = '  +2.94%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.132'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +15.59%  '
